# Apply update to the "Registros" sheet:
#  - Row 9 (Joao Pedro Santos): update several columns to new values
#  - Row 13 (new): add Rodrigo Souza Santos record

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

function Set-TextCell {
    param($sheet, [string]$addr, [string]$val)
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Update existing row 9 ---
Set-TextCell $ws "B9" "0"
Set-TextCell $ws "C9" "preto"
Set-TextCell $ws "D9" "0"
Set-TextCell $ws "E9" "0"
Set-TextCell $ws "F9" "0"
Set-TextCell $ws "G9" "0"
Set-TextCell $ws "I9" "0"

# --- Add new row 13 ---
Set-TextCell $ws "A13" "Rodrigo Souza Santos"
Set-TextCell $ws "B13" "26"
Set-TextCell $ws "C13" "no"
Set-TextCell $ws "D13" "1"
Set-TextCell $ws "E13" "1"
Set-TextCell $ws "F13" "1"
Set-TextCell $ws "G13" "0"
Set-TextCell $ws "H13" "0"
Set-TextCell $ws "I13" "0"
Set-TextCell $ws "J13" "0"
Set-TextCell $ws "K13" "1"
Set-TextCell $ws "L13" "0"
Set-TextCell $ws "M13" "0"
Set-TextCell $ws "N13" "0"
Set-TextCell $ws "O13" "0"
Set-TextCell $ws "P13" "0"
Set-TextCell $ws "Q13" "0"
Set-TextCell $ws "R13" "0"
Set-TextCell $ws "S13" "0"
Set-TextCell $ws "T13" "0"
Set-TextCell $ws "U13" "0"
Set-TextCell $ws "V13" "1"
Set-TextCell $ws "W13" "0"
